$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of new-string introduction matches target sharedStrings.xml ordering:
# Phobos, distToParent [10^9 m], 10.6*10^15, Deimos, 1.4762x10^15, Moon

# 4. Row 14: Phobos (name introduced first)
$ws.Range("A14").Value = "Phobos"

# 1. Rename header G1 from "distToSun [10^9 m]" to "distToParent [10^9 m]"
$ws.Range("G1").Value = "distToParent [10^9 m]"

$ws.Range("B14").Value = "10.6*10^15"
$ws.Range("D14").Value = 5
$ws.Range("G14").Value = 9375000
$ws.Range("H14").Value = 0.32

# 5. Row 15: Deimos data
$ws.Range("A15").Value = "Deimos"
$ws.Range("B15").Value = "1.4762×10^15"
$ws.Range("D15").Value = 6

# 3. Row 13: Moon data (name introduced last)
$ws.Range("A13").Value = "Moon"
$ws.Range("B13").Value = 0.07346
$ws.Range("D13").Value = 1738.1
$ws.Range("E13").Value = 0.3633
$ws.Range("F13").Value = 0.4055
$ws.Range("G13").Formula = "=(E13+F13)/2"
$ws.Range("H13").Value = 27.3217
$ws.Range("E13").NumberFormat = "#,##0.00"
$ws.Range("F13").NumberFormat = "#,##0.00"
$ws.Range("H13").NumberFormat = "#,##0"

# 2. Row 12 acts as a blank separator row but carries number formatting
$ws.Range("E12").NumberFormat = "#,##0.00"
$ws.Range("F12").NumberFormat = "#,##0.00"
$ws.Range("H12").NumberFormat = "#,##0"

# 6. Update selection to G13
$ws.Range("G13").Select() | Out-Null

Write-Output "done"
